$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spain Primera Liga")

function Swap-Row($r1, $r2) {
    $cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Swap rows 169 and 170
Swap-Row 169 170

# Swap rows 189 and 190
Swap-Row 189 190

# Row 311 individual cell updates
$ws.Range("N311").Value = 3.3
$ws.Range("P311").Value = 2.375
$ws.Range("R311").Value = 1.84
$ws.Range("S311").Value = 2.06
